# Append a new row (row 96) of data to each of the four worksheets,
# mirroring the structure of the existing row 95 on each sheet.
#
# Per-sheet values for the new row:
#   Sheet 1 (MID_LFT_#1): B=0x01,0x90  C=0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,  D=0x01,0x18  E=0x07  F=400  G=5.68631262647113e+23  H=280  I=7
#   Sheet 2 (MID_LFT_#2): B=0x01,0x7c  C=0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,  D=0x01,0x24  E=0x19  F=380  G=5.68432987514711e+23  H=292  I=25
#   Sheet 3 (MID_PLT_#1): B=0x00,0x6e  C=0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,  D=0x00,0x5E  E=0x15  F=110  G=5.68631262647113e+23  H=94   I=15
#   Sheet 4 (MID_PLT_#2): B=0x00,0x82  C=0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,  D=0x00,0x74  E=0x9   F=130  G=5.68631262647113e+23  H=116  I=9
#
# All four sheets share the same new timestamp for column A: 45882.46402777778
# (formatted the same way as the date already used for column A elsewhere).

$wb = $excel.ActiveWorkbook

$newRow = 96
$prevRow = 95
$newDate = 45882.46402777778

$rowData = @{
    1 = @{ B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; D = "0x01,0x18"; E = "0x07";  F = 400; G = "5.68631262647113e+23"; H = 280; I = 7 }
    2 = @{ B = "0x01,0x7c"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; D = "0x01,0x24"; E = "0x19";  F = 380; G = "5.68432987514711e+23"; H = 292; I = 25 }
    3 = @{ B = "0x00,0x6e"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; D = "0x00,0x5E"; E = "0x15";  F = 110; G = "5.68631262647113e+23"; H = 94;  I = 15 }
    4 = @{ B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; D = "0x00,0x74"; E = "0x9";   F = 130; G = "5.68631262647113e+23"; H = 116; I = 9 }
}

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $data = $rowData[$i]

    # Column A: same date/time number format as the rest of column A.
    $ws.Cells.Item($newRow, 1).Value = $newDate
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat

    # Columns B-E: text values (hex byte lists).
    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E

    # Columns F-I: numeric values.
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = [double]$data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
